# Updated README to include Wish List
# Adds a new "Wish List Track List" worksheet after "Job Status", containing
# a wish-list of tracks (mirroring the Tracks sheet layout/formatting), and
# updates the Tracks sheet selection to a full-column selection.

$wb = $excel.ActiveWorkbook

$tracksSheet = $wb.Worksheets.Item("Tracks")
$jobStatusSheet = $wb.Worksheets.Item("Job Status")

# --- 1. Add the new worksheet after "Job Status" ---------------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $jobStatusSheet)
$ws.Name = "Wish List Track List"

# --- 2. Seed layout/formatting by copying the Tracks sheet's header +      --
#        12 data rows (A1:E13), which carries over the same cell styles    --
#        (bold header row, left-aligned text, left-aligned time format).   --
$tracksSheet.Range("A1:E13").Copy($ws.Range("A1"))

# Also stamp styled (but empty) duration cells for rows 14-18, matching the
# "leftover formatted rows" below the data in the source file.
$tracksSheet.Range("E2:E6").Copy($ws.Range("E14"))
$ws.Range("E14:E18").ClearContents()

# --- 3. Overwrite the copied data with the real Wish List track data ------
$wishListData = @(
    @("Live In Paris", "Diana Krall", 1,  "I Love Being Here With You", 0.21666666666666667),
    @("Live In Paris", "Diana Krall", 2,  "Let's Fall In Love", 0.19027777777777777),
    @("Live In Paris", "Diana Krall", 3,  "'Deed I Do", 0.22013888888888888),
    @("Live In Paris", "Diana Krall", 4,  "The Look Of Love", 0.20833333333333334),
    @("Live In Paris", "Diana Krall", 5,  "East Of The Sun (And West Of The Moon)", 0.24861111111111112),
    @("Live In Paris", "Diana Krall", 6,  "I've Got You Under My Skin", 0.30833333333333335),
    @("Live In Paris", "Diana Krall", 7,  "Devil May Care", 0.28611111111111115),
    @("Live In Paris", "Diana Krall", 8,  "Maybe You'll Be There", 0.24097222222222223),
    @("Live In Paris", "Diana Krall", 9,  "'S Wonderful", 0.24930555555555556),
    @("Live In Paris", "Diana Krall", 10, "Fly Me To The Moon", 0.25347222222222221),
    @("Live In Paris", "Diana Krall", 11, "A Case Of You", 0.29444444444444445),
    @("Live In Paris", "Diana Krall", 12, "Just The Way You Are", 0.20833333333333334)
)

$row = 2
foreach ($track in $wishListData) {
    $ws.Range("A$row").Value = $track[0]
    $ws.Range("B$row").Value = $track[1]
    $ws.Range("C$row").Value = $track[2]
    $ws.Range("D$row").Value = $track[3]
    $ws.Range("E$row").Value = $track[4]
    $row++
}

# --- 4. Match column widths to the Tracks sheet (best effort) -------------
$ws.Columns.Item(1).ColumnWidth = $tracksSheet.Columns.Item(1).ColumnWidth
$ws.Columns.Item(2).ColumnWidth = $tracksSheet.Columns.Item(2).ColumnWidth
$ws.Columns.Item(3).ColumnWidth = $tracksSheet.Columns.Item(3).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $tracksSheet.Columns.Item(4).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $tracksSheet.Columns.Item(5).ColumnWidth

# --- 5. Update the Tracks sheet selection to a full-column selection ------
$tracksSheet.Range("A1:XFD1048576").Select()

# --- 6. Leave the new sheet as the active tab ------------------------------
$ws.Select()
